# Liquidación HFF_Liquidation-8F_Air-AWB_157-84730612
#
# Fill in the missing "Weight" (重量) value for every pallet row of the
# main control table. Row 10 already carries the box weight "2.5" (kg) for
# this lot; rows 11-18 (the remaining pallets of lot 157-84730612,
# including the "Custom"/no vendidos ones) were left blank in column E.
# Complete them with the same quantitative value so the whole block is
# consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

for ($row = 11; $row -le 18; $row++) {
    $ws.Cells.Item($row, 5).Value = "2.5"
}
